$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (columns A-G)
$data = @(
    @(34525, "Maria Eduarda Gomes", "Operacoes", "Viagem de negocios", 7, 45095, 9865.280000000001),
    @(13446, "Aurora Martins", "Financeiro", "Outros", 2, 45103, 4636.8),
    @(71714, "Srta. Laís Cardoso", "TI", "Problemas pessoais", 7, 45086, 7281.97),
    @(35287, "Melissa Lima", "Financeiro", "Doenca", 7, 45097, 4037.71),
    @(33665, "Benjamin Moura", "Atendimento ao Cliente", "Viagem de negocios", 5, 45090, 2970.79),
    @(56083, "Maria Luiza Barros", "Atendimento ao Cliente", "Outros", 8, 45087, 6260.23),
    @(57896, "Vinicius Abreu", "Financeiro", "Problemas pessoais", 5, 45079, 4412.39),
    @(90525, "Calebe Lopes", "Recursos Humanos", "Outros", 8, 45094, 8757.139999999999),
    @(1803, "Sra. Isabella Sales", "Vendas", "Outros", 6, 45083, 2440.77),
    @(70185, "Srta. Heloísa Marques", "TI", "Problemas pessoais", 1, 45086, 6965.41)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
}
